$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was updated from
# 45175 (2023-09-06) to 45183 (2023-09-14) for every data row (2-158).
for ($r = 2; $r -le 158; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45175) {
        $cell.Value = 45183
    }
}
